$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "The user gets if the website." -> "The user gets on the website."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("The user gets if the website.", $false, $false, $false, $false, $false, $true, 1, $false, "The user gets on the website.", 2) | Out-Null

Write-Host "P3:" $d.Paragraphs.Item(3).Range.Text

# ---------------------------------------------------------------------------
# 2) Insert three new paragraphs after "The user has the option to play
#    preloaded contents that are made by me." (paragraph 5), before
#    "The user scrolls down..." (paragraph 6)
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5).Range
$p5.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(6).Range
$newPara.Text = "The user can share content that is uploaded b me to social networks. "

$newPara.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item(7).Range
$newPara2.Text = "The user can download songs from the website."

$newPara2.InsertParagraphAfter()
$newPara3 = $d.Paragraphs.Item(8).Range
$newPara3.Text = "User can post comments on the page."

Write-Host "Count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  Write-Host $i ": [" $d.Paragraphs.Item($i).Range.Text "]"
}

# ---------------------------------------------------------------------------
# 3) "The user can share any song on the page." -> "The user can like any
#    song on the page."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("The user can share any song on the page.", $false, $false, $false, $false, $false, $true, 1, $false, "The user can like any song on the page.", 2) | Out-Null

$sharePara = $d.Paragraphs.Item(10)
Write-Host "P10:" $sharePara.Range.Text

# ---------------------------------------------------------------------------
# 4) Insert new paragraph after it: "This website will bring together fans
#    of Lo-fi music. "
# ---------------------------------------------------------------------------
$sharePara.Range.InsertParagraphAfter()
$newPara4 = $d.Paragraphs.Item(11).Range
$newPara4.Text = "This website will bring together fans of Lo-fi music. "

Write-Host "Count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  Write-Host $i ": [" $d.Paragraphs.Item($i).Range.Text "]"
}
